$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-06-06 Thursday" "2024-06-07 Friday"

Replace-Text "307×3=" "243×2="
Replace-Text "279×8=" "443×9="
Replace-Text "915×6=" "821×8="
Replace-Text "267×7=" "956×5="
Replace-Text "364×8=" "495×8="
Replace-Text "943×9=" "231×4="
Replace-Text "757×3=" "888×5="
Replace-Text "651×3=" "739×4="
Replace-Text "803×2=" "498×7="
Replace-Text "176×7=" "826×8="
Replace-Text "855×7=" "898×6="
Replace-Text "648×8=" "678×3="
Replace-Text "795×7=" "664×4="
Replace-Text "275×4=" "881×9="
Replace-Text "840×2=" "765×4="
Replace-Text "525×2=" "139×2="
Replace-Text "436×8=" "360×8="
Replace-Text "990×9=" "475×2="
Replace-Text "822×6=" "847×3="
Replace-Text "611×2=" "318×6="
Replace-Text "222×5=" "438×8="
Replace-Text "796×6=" "631×6="
Replace-Text "804×5=" "332×7="
Replace-Text "803×6=" "198×2="
Replace-Text "836×4=" "158×4="
